$wb = $excel.ActiveWorkbook

# "Metadata" sheet updates
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/er-or-observation-room-service"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# "Include from Er Or Observatio" sheet updates
$codes = $wb.Worksheets.Item("Include from Er Or Observatio")
$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/er-or-observation-room-service"
